# Daily auto push update: insert two new rows of data for 2026/02/25 (水)
# right after the existing 2026/02/25 rows (old row 856), shifting all the
# following rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 857; this pushes the former
# rows 857-898 down to 859-900.
$ws.Rows("857:858").Insert()

# --- New row 857 ---
$ws.Range("A857").NumberFormat = "@"
$ws.Range("A857").Value = "2026/02/25"
$ws.Range("A857").Style = "Normal"
$ws.Range("B857").Value = "水"
$ws.Range("C857").Value = 19
$ws.Range("D857").Value = 39

# --- New row 858 ---
$ws.Range("A858").NumberFormat = "@"
$ws.Range("A858").Value = "2026/02/25"
$ws.Range("A858").Style = "Normal"
$ws.Range("B858").Value = "水"
$ws.Range("C858").Value = 22
$ws.Range("D858").Value = 42
